$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (index 3) and column E (index 5) values per the diff
$values = @{
    2 = @{ C = 3.123541145015474; E = 4.320516327661528 }
    3 = @{ C = 3.959010658874851; E = 3.433494243648449 }
    4 = @{ C = 4.722695063536686; E = 7.855477094481422 }
    5 = @{ C = 6.739021039846627; E = -0.7259153295281151 }
    6 = @{ C = 2.619839412265601; E = -0.17790865651377 }
    7 = @{ C = -0.7919564768266385; E = 2.233697987812078 }
    8 = @{ C = 1.877689851450803; E = 2.959667200710037 }
    9 = @{ C = 2.705004599189187; E = 0.5406927319912658 }
    10 = @{ C = 1.110374544249249; E = 3.29487077883559 }
    11 = @{ C = 2.267566233338814; E = 3.15890982365572 }
    12 = @{ C = 2.688433258834588; E = 3.239674285955152 }
    13 = @{ C = 1.014079695989589; E = 2.571626871154176 }
    14 = @{ C = 3.013853578092252; E = 2.652928973511215 }
    15 = @{ C = 1.331295149770684; E = -0.3738725857433511 }
    16 = @{ C = 0.04589006555719699; E = 1.021287096146906 }
    17 = @{ C = 0.009546395482029624; E = 0.5726247744375135 }
    18 = @{ C = 0.8709390141433015; E = 0.6379602509701376 }
    19 = @{ C = 0.7652063367885598; E = 2.297544413125596 }
    20 = @{ C = 2.267579219134386; E = 2.825914290412324 }
    21 = @{ C = 3.146753122914103; E = 2.510325059131513 }
    22 = @{ C = 1.769033835366818; E = -1.563640406432543 }
    23 = @{ C = -4.774715709990263; E = 1.685921024959058 }
    24 = @{ C = 1.95493704440024; E = 3.595026567604331 }
    25 = @{ C = 3.478075069442799; E = 2.123512403013161 }
    26 = @{ C = 1.232342134690434; E = 1.023960954496861 }
    27 = @{ C = 0.2542811494408159; E = 1.414810393331356 }
    28 = @{ C = 1.519778766382096; E = 0.3626330124320232 }
    29 = @{ C = 1.469441753880329; E = 1.339091979913909 }
    30 = @{ C = 1.638203081492495; E = 1.006355688239569 }
    31 = @{ C = 2.268697431234346; E = 3.257368055312471 }
    32 = @{ C = 1.984425467899631; E = -0.3299132127116078 }
    33 = @{ C = 0.6066448776129052; E = 0.2691274977562275 }
    34 = @{ C = -4.243076347305386; E = 23.52713615747899 }
    35 = @{ C = 1.438499295329754; E = 6.182015844361843 }
    36 = @{ C = 1.906593537051537; E = 1.559766133975371 }
    37 = @{ C = 0.08348019664223827; E = -0.3513378667146627 }
    38 = @{ C = -0.214505326882275; E = 0.256124515548195 }
    39 = @{ C = 0.1651547428133782; E = -0.02867520550564606 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row].C
    $ws.Cells.Item($row, 5).Value = $values[$row].E
}

$wb.Save()